# Insert three new indicator rows (INDICATOR_81, INDICATOR_82, INDICATOR_83)
# at the top of the INDICATOR_200.. block on the "Library_Formula" sheet,
# pushing the existing rows (formerly 84-159) down to 87-162.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert 3 blank rows before row 84 (Excel copies formatting from the row above).
$ws.Rows("84:86").Insert()

$newIndicators = @("INDICATOR_81", "INDICATOR_82", "INDICATOR_83")

for ($i = 0; $i -lt $newIndicators.Count; $i++) {
    $r = 84 + $i
    $ws.Cells.Item($r, 1).Value = "CREATE/MODIFY"
    $ws.Cells.Item($r, 2).Value = "LIB_EWS_IT"
    $ws.Cells.Item($r, 3).Value = $newIndicators[$i]
    $ws.Cells.Item($r, 5).Value = "String"
    $ws.Cells.Item($r, 6).Value = "String"
}

# Match the saved view state from the diff (selection near the new rows).
$ws.Activate()
$ws.Range("E83:F86").Select()
